$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "67.764.26"

Set-TextValue "D3" "3.805.68"
Set-TextValue "E3" "  +0.39%  "

Set-TextValue "E4" "  +0.08%  "

Set-TextValue "D5" "607.73"
Set-TextValue "E5" "  +2.10%  "

Set-TextValue "D6" "166.90"
Set-TextValue "E6" "  +0.26%  "

Set-TextValue "E7" "  +0.14%  "

Set-TextValue "E8" "  +0.48%  "

Set-TextValue "E9" "  +0.89%  "

Set-TextValue "E10" "  -0.75%  "

Set-TextValue "E11" "  +0.13%  "

Set-TextValue "E12" "  -0.97%  "

Set-TextValue "D13" "35.96"
Set-TextValue "E13" "  -0.84%  "

Set-TextValue "D14" "4.443.89"
Set-TextValue "E14" "  +0.29%  "

Set-TextValue "D15" "3.825.03"
Set-TextValue "E15" "  +0.92%  "

Set-TextValue "D16" "18.49"
Set-TextValue "E16" "  -0.03%  "

Set-TextValue "D17" "67.860.55"
Set-TextValue "E17" "  +0.33%  "

Set-TextValue "D18" "7.09"
Set-TextValue "E18" "  +1.41%  "

Set-TextValue "E19" "  +0.48%  "

Set-TextValue "D20" "462.01"
Set-TextValue "E20" "  +1.04%  "

Set-TextValue "D21" "9.90"
Set-TextValue "E21" "  -3.16%  "

Set-TextValue "D22" "0.701"
Set-TextValue "E22" "  +0.59%  "

Set-TextValue "E23" "  -2.45%  "

Set-TextValue "D24" "83.37"
Set-TextValue "E24" "  -0.11%  "

Set-TextValue "D25" "12.09"
Set-TextValue "E25" "  +1.31%  "

Set-TextValue "D26" "2.11"
Set-TextValue "E26" "  -1.47%  "

Set-TextValue "E27" "  +0.13%  "

Set-TextValue "D28" "10.01"
Set-TextValue "E28" "  -0.80%  "

Set-TextValue "D29" "3.955.21"
Set-TextValue "E29" "  +0.35%  "

Set-TextValue "D30" "2.80"
Set-TextValue "E30" "  +0.46%  "

Set-TextValue "D31" "7.40"
Set-TextValue "E31" "  +1.44%  "

Set-TextValue "D32" "2.23"
Set-TextValue "E32" "  +1.65%  "

Set-TextValue "D33" "29.58"
Set-TextValue "E33" "  -0.75%  "

Set-TextValue "E35" "  -1.44%  "

Set-TextValue "D36" "3.746.44"
Set-TextValue "E36" "  +0.02%  "

Set-TextValue "D37" "0.100"
Set-TextValue "E37" "  +0.00%  "

Set-TextValue "D38" "3.38"
Set-TextValue "E38" "  +1.45%  "

Set-TextValue "E39" "  -0.24%  "

Set-TextValue "E40" "  +0.44%  "

Set-TextValue "E41" "  +0.51%  "

Set-TextValue "D42" "1.00"
Set-TextValue "E42" "  +0.02%  "

Set-TextValue "E43" "  -0.01%  "

Set-TextValue "D44" "48.19"
Set-TextValue "E44" "  +2.25%  "

Set-TextValue "E45" "  +0.83%  "

Set-TextValue "D46" "43.05"
Set-TextValue "E46" "  -4.14%  "

Set-TextValue "D47" "28.00"
Set-TextValue "E47" "  +10.61%  "

Set-TextValue "D48" "8.35"
Set-TextValue "E48" "  -0.31%  "

Set-TextValue "E49" "  +9.54%  "

Set-TextValue "D50" "148.57"
Set-TextValue "E50" "  -0.17%  "

Set-TextValue "E51" "  +0.39%  "

